# Auto-generated edit script: apply scraped-data update to horarios workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$rows = @(
    @(1, "LÍNEA 141 - LP1912 - 16/01/2026", $null, $null, $null, $null),
    @(2, "Última actualización: 08:33:47", $null, $null, $null, $null),
    @(3, "Total filas: 121", $null, $null, $null, $null),
    @(5, "Hora_Scrap", "Hora_Llegada", "Linea", "Minutos", "Parada"),
    @(6, "00:04:05", "01:12", "215_ALUAR", 68, "LP1912"),
    @(7, "00:04:05", "01:58", "14_ABASTO", 114, "LP1912"),
    @(8, "01:07:17", "02:58", "215_ALUAR", 111, "LP1912"),
    @(9, "01:49:17", "03:03", "215_ALUAR", 74, "LP1912"),
    @(10, "01:49:17", "03:48", "14_ABASTO", 119, "LP1912"),
    @(11, "02:42:39", "03:54", "14_ABASTO", 72, "LP1912"),
    @(12, "02:13:28", "03:55", "14_ABASTO", 102, "LP1912"),
    @(13, "02:13:28", "04:01", "81_EL PELIGRO", 108, "LP1912"),
    @(14, "03:23:38", "04:45", "215A_EL PATO", 82, "LP1912"),
    @(15, "02:56:05", "04:46", "215A_EL PATO", 110, "LP1912"),
    @(16, "02:56:05", "04:53", "11_ETCHEVERRY", 117, "LP1912"),
    @(17, "04:53:50", "04:54", "11_ETCHEVERRY", 1, "LP1912"),
    @(18, "04:53:50", "05:15", "14_ABASTO", 22, "LP1912"),
    @(19, "03:23:38", "05:16", "17_ROMERO", 113, "LP1912"),
    @(20, "04:40:32", "05:17", "14_ABASTO", 37, "LP1912"),
    @(21, "05:18:23", "05:20", "14_ABASTO", 2, "LP1912"),
    @(22, "03:23:38", "05:22", "23_HERNANDEZ", 119, "LP1912"),
    @(23, "03:54:40", "05:34", "215B_EL PATO", 100, "LP1912"),
    @(24, "04:18:06", "05:35", "215B_EL PATO", 77, "LP1912"),
    @(25, "04:18:06", "05:36", "14_ABASTO", 78, "LP1912"),
    @(26, "03:54:40", "05:39", "14_ABASTO", 105, "LP1912"),
    @(27, "03:54:40", "05:46", "15_ABASTO", 112, "LP1912"),
    @(28, "05:47:32", "05:47", "15_ABASTO", 0, "LP1912"),
    @(29, "05:47:32", "05:49", "14_ABASTO", 2, "LP1912"),
    @(30, "04:40:32", "06:04", "16_SANTA ANA", 84, "LP1912"),
    @(31, "06:02:16", "06:05", "16_SANTA ANA", 3, "LP1912"),
    @(32, "04:18:06", "06:09", "16_SANTA ANA", 111, "LP1912"),
    @(33, "04:40:32", "06:11", "215A_EL PATO", 91, "LP1912"),
    @(34, "04:18:06", "06:12", "215A_EL PATO", 114, "LP1912"),
    @(35, "04:18:06", "06:14", "225_HARAS DEL SUR", 116, "LP1912"),
    @(36, "04:40:32", "06:21", "26_HERNANDEZ", 101, "LP1912"),
    @(37, "04:40:32", "06:27", "23_HERNANDEZ", 107, "LP1912"),
    @(38, "04:40:32", "06:29", "86_EST CHICA-ESC AGRARIA", 109, "LP1912"),
    @(39, "06:02:16", "06:30", "86_EST CHICA-ESC AGRARIA", 28, "LP1912"),
    @(40, "04:40:32", "06:31", "16_SANTA ANA", 111, "LP1912"),
    @(41, "04:53:50", "06:44", "225_C ROCA-H SUR", 111, "LP1912"),
    @(42, "04:53:50", "06:46", "215C_EL PATO", 113, "LP1912"),
    @(43, "05:18:23", "06:58", "10_OLMOS", 100, "LP1912"),
    @(44, "05:18:23", "06:59", "14_ABASTO", 101, "LP1912"),
    @(45, "06:02:16", "07:00", "14_ABASTO", 58, "LP1912"),
    @(46, "06:37:24", "07:01", "16_SANTA ANA", 24, "LP1912"),
    @(47, "05:47:32", "07:04", "23_HERNANDEZ", 77, "LP1912"),
    @(48, "06:02:16", "07:05", "23_HERNANDEZ", 63, "LP1912"),
    @(49, "05:18:23", "07:05", "15_ABASTO", 107, "LP1912"),
    @(50, "05:18:23", "07:07", "225_GOMEZ", 109, "LP1912"),
    @(51, "05:18:23", "07:11", "215A_EL PATO", 113, "LP1912"),
    @(52, "07:14:27", "07:14", "11_ETCHEVERRY", 0, "LP1912"),
    @(53, "05:18:23", "07:15", "11_ETCHEVERRY", 117, "LP1912"),
    @(54, "06:02:16", "07:16", "11_ETCHEVERRY", 74, "LP1912"),
    @(55, "06:37:24", "07:16", "16_SANTA ANA", 39, "LP1912"),
    @(56, "05:47:32", "07:21", "26_HERNANDEZ", 94, "LP1912"),
    @(57, "06:02:16", "07:23", "10_OLMOS", 81, "LP1912"),
    @(58, "05:47:32", "07:27", "10_OLMOS", 100, "LP1912"),
    @(59, "05:47:32", "07:31", "16_SANTA ANA", 104, "LP1912"),
    @(60, "05:47:32", "07:31", "11_ETCHEVERRY", 104, "LP1912"),
    @(61, "06:02:16", "07:32", "11_ETCHEVERRY", 90, "LP1912"),
    @(62, "05:47:32", "07:32", "84_COLONIA URQUIZA-ESC 49", 105, "LP1912"),
    @(63, "06:37:24", "07:34", "23_HERNANDEZ", 57, "LP1912"),
    @(64, "05:47:32", "07:36", "27_EL RETIRO", 109, "LP1912"),
    @(65, "06:02:16", "07:37", "27_EL RETIRO", 95, "LP1912"),
    @(66, "07:14:27", "07:37", "23_HERNANDEZ", 23, "LP1912"),
    @(67, "05:47:32", "07:39", "10_OLMOS", 112, "LP1912"),
    @(68, "07:44:08", "07:45", "14_ABASTO", 1, "LP1912"),
    @(69, "06:37:24", "07:47", "14_ABASTO", 70, "LP1912"),
    @(70, "06:02:16", "07:48", "14_ABASTO", 106, "LP1912"),
    @(71, "06:02:16", "07:51", "215D_EL PATO", 109, "LP1912"),
    @(72, "07:44:08", "07:52", "215D_EL PATO", 8, "LP1912"),
    @(73, "07:44:08", "07:55", "10_OLMOS", 11, "LP1912"),
    @(74, "07:14:27", "07:58", "16_SANTA ANA", 44, "LP1912"),
    @(75, "07:14:27", "08:03", "11_ETCHEVERRY", 49, "LP1912"),
    @(76, "07:44:08", "08:04", "11_ETCHEVERRY", 20, "LP1912"),
    @(77, "07:57:27", "08:06", "11_ETCHEVERRY", 9, "LP1912"),
    @(78, "07:57:27", "08:10", "16_SANTA ANA", 13, "LP1912"),
    @(79, "07:44:08", "08:11", "16_SANTA ANA", 27, "LP1912"),
    @(80, "06:37:24", "08:12", "15_ABASTO", 95, "LP1912"),
    @(81, "07:44:08", "08:13", "10_OLMOS", 29, "LP1912"),
    @(82, "06:37:24", "08:21", "26_HERNANDEZ", 104, "LP1912"),
    @(83, "06:37:24", "08:22", "16_P MOR-SANTA ANA", 105, "LP1912"),
    @(84, "06:37:24", "08:23", "215B_EL PATO", 106, "LP1912"),
    @(85, "07:44:08", "08:23", "16_P MOR-SANTA ANA", 39, "LP1912"),
    @(86, "06:37:24", "08:27", "84_COLONIA URQUIZA-ESC 49", 110, "LP1912"),
    @(87, "07:57:27", "08:33", "10_OLMOS", 36, "LP1912"),
    @(88, "06:52:38", "08:42", "81_EL PELIGRO", 110, "LP1912"),
    @(89, "07:14:27", "08:43", "14_ABASTO", 89, "LP1912"),
    @(90, "07:44:08", "08:44", "14_ABASTO", 60, "LP1912"),
    @(91, "07:14:27", "08:54", "17_ROMERO", 100, "LP1912"),
    @(92, "08:16:48", "08:55", "10_OLMOS", 39, "LP1912"),
    @(93, "07:14:27", "09:01", "215A_EL PATO", 107, "LP1912"),
    @(94, "07:44:08", "09:02", "215A_EL PATO", 78, "LP1912"),
    @(95, "07:57:27", "09:03", "11_ETCHEVERRY", 66, "LP1912"),
    @(96, "08:16:48", "09:04", "23_HERNANDEZ", 48, "LP1912"),
    @(97, "07:14:27", "09:07", "23_HERNANDEZ", 113, "LP1912"),
    @(98, "07:44:08", "09:08", "23_HERNANDEZ", 84, "LP1912"),
    @(99, "07:57:27", "09:09", "23_HERNANDEZ", 72, "LP1912"),
    @(100, "07:14:27", "09:10", "16_P MOR-SANTA ANA", 116, "LP1912"),
    @(101, "07:44:08", "09:11", "16_P MOR-SANTA ANA", 87, "LP1912"),
    @(102, "08:33:47", "09:13", "10_OLMOS", 40, "LP1912"),
    @(103, "07:44:08", "09:14", "16_SANTA ANA", 90, "LP1912"),
    @(104, "07:57:27", "09:16", "27_EL RETIRO", 79, "LP1912"),
    @(105, "07:44:08", "09:17", "27_EL RETIRO", 93, "LP1912"),
    @(106, "07:44:08", "09:21", "26_HERNANDEZ", 97, "LP1912"),
    @(107, "07:57:27", "09:22", "16_SANTA ANA", 85, "LP1912"),
    @(108, "07:57:27", "09:22", "17_ROMERO", 85, "LP1912"),
    @(109, "07:44:08", "09:23", "17_ROMERO", 99, "LP1912"),
    @(110, "07:57:27", "09:23", "11_ETCHEVERRY", 86, "LP1912"),
    @(111, "07:44:08", "09:24", "11_ETCHEVERRY", 100, "LP1912"),
    @(112, "08:16:48", "09:29", "16_SANTA ANA", 73, "LP1912"),
    @(113, "07:44:08", "09:32", "15_ABASTO", 108, "LP1912"),
    @(114, "07:44:08", "09:33", "10_OLMOS", 109, "LP1912"),
    @(115, "08:33:47", "09:34", "16_SANTA ANA", 61, "LP1912"),
    @(116, "07:44:08", "09:36", "23_HERNANDEZ", 112, "LP1912"),
    @(117, "08:16:48", "09:37", "23_HERNANDEZ", 81, "LP1912"),
    @(118, "08:16:48", "09:41", "215C_EL PATO", 85, "LP1912"),
    @(119, "08:33:47", "09:41", "23_HERNANDEZ", 68, "LP1912"),
    @(120, "07:44:08", "09:42", "215C_EL PATO", 118, "LP1912"),
    @(121, "07:57:27", "09:43", "14_ABASTO", 106, "LP1912"),
    @(122, "08:16:48", "10:10", "16_P MOR-SANTA ANA", 114, "LP1912"),
    @(123, "08:16:48", "10:12", "15_ABASTO", 116, "LP1912"),
    @(124, "08:33:47", "10:21", "26_HERNANDEZ", 108, "LP1912"),
    @(125, "08:33:47", "10:22", "17_ROMERO", 109, "LP1912"),
    @(126, "08:33:47", "10:26", "215A_EL PATO", 113, "LP1912")
)
foreach ($row in $rows) {
    $r = $row[0]
    if ($row[1] -ne $null) { $ws.Cells.Item($r, 1).Value = $row[1] }
    if ($row[2] -ne $null) { $ws.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 3).Value = $row[3] }
    if ($row[4] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[4] }
    if ($row[5] -ne $null) { $ws.Cells.Item($r, 5).Value = $row[5] }
}

$ws = $wb.Worksheets.Item("LP1912-215")
$rows = @(
    @(1, "LÍNEA 141 - LP1912-215 - 16/01/2026", $null, $null, $null, $null),
    @(2, "Última actualización: 08:33:47", $null, $null, $null, $null),
    @(3, "Total filas: 19", $null, $null, $null, $null),
    @(5, "Hora_Scrap", "Hora_Llegada", "Linea", "Minutos", "Parada"),
    @(6, "00:04:05", "01:12", "215_ALUAR", 68, "LP1912"),
    @(7, "01:07:17", "02:58", "215_ALUAR", 111, "LP1912"),
    @(8, "01:49:17", "03:03", "215_ALUAR", 74, "LP1912"),
    @(9, "03:23:38", "04:45", "215A_EL PATO", 82, "LP1912"),
    @(10, "02:56:05", "04:46", "215A_EL PATO", 110, "LP1912"),
    @(11, "03:54:40", "05:34", "215B_EL PATO", 100, "LP1912"),
    @(12, "04:18:06", "05:35", "215B_EL PATO", 77, "LP1912"),
    @(13, "04:40:32", "06:11", "215A_EL PATO", 91, "LP1912"),
    @(14, "04:18:06", "06:12", "215A_EL PATO", 114, "LP1912"),
    @(15, "04:53:50", "06:46", "215C_EL PATO", 113, "LP1912"),
    @(16, "05:18:23", "07:11", "215A_EL PATO", 113, "LP1912"),
    @(17, "06:02:16", "07:51", "215D_EL PATO", 109, "LP1912"),
    @(18, "07:44:08", "07:52", "215D_EL PATO", 8, "LP1912"),
    @(19, "06:37:24", "08:23", "215B_EL PATO", 106, "LP1912"),
    @(20, "07:14:27", "09:01", "215A_EL PATO", 107, "LP1912"),
    @(21, "07:44:08", "09:02", "215A_EL PATO", 78, "LP1912"),
    @(22, "08:16:48", "09:41", "215C_EL PATO", 85, "LP1912"),
    @(23, "07:44:08", "09:42", "215C_EL PATO", 118, "LP1912"),
    @(24, "08:33:47", "10:26", "215A_EL PATO", 113, "LP1912")
)
foreach ($row in $rows) {
    $r = $row[0]
    if ($row[1] -ne $null) { $ws.Cells.Item($r, 1).Value = $row[1] }
    if ($row[2] -ne $null) { $ws.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 3).Value = $row[3] }
    if ($row[4] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[4] }
    if ($row[5] -ne $null) { $ws.Cells.Item($r, 5).Value = $row[5] }
}

$ws = $wb.Worksheets.Item("6203-6173")
$rows = @(
    @(1, "LÍNEA 141 - 6203-6173 - 16/01/2026", $null, $null, $null, $null),
    @(2, "Última actualización: 08:33:47", $null, $null, $null, $null),
    @(3, "Total filas: 22", $null, $null, $null, $null),
    @(5, "Hora_Scrap", "Hora_Llegada", "Linea", "Minutos", "Parada"),
    @(6, "00:04:05", "00:05", "215A_LA PLATA", 1, "L6173"),
    @(7, "03:54:40", "05:43", "215A_LA PLATA", 109, "L6173"),
    @(8, "04:18:06", "05:44", "215A_LA PLATA", 86, "L6173"),
    @(9, "04:53:50", "06:08", "215A_LA PLATA", 75, "L6173"),
    @(10, "04:18:06", "06:09", "215A_LA PLATA", 111, "L6173"),
    @(11, "04:53:50", "06:32", "215C_LA PLATA", 99, "L6203"),
    @(12, "04:40:32", "06:33", "215C_LA PLATA", 113, "L6203"),
    @(13, "05:18:23", "07:00", "215B_LP-P MOR-1 Y 57", 102, "L6173"),
    @(14, "05:47:32", "07:35", "215A_LA PLATA", 108, "L6173"),
    @(15, "06:52:38", "08:06", "215C_LA PLATA", 74, "L6203"),
    @(16, "06:37:24", "08:07", "215C_LA PLATA", 90, "L6203"),
    @(17, "07:14:27", "08:14", "215C_LA PLATA", 60, "L6203"),
    @(18, "07:44:08", "08:17", "215C_LA PLATA", 33, "L6203"),
    @(19, "07:57:27", "08:18", "215C_LA PLATA", 21, "L6203"),
    @(20, "06:37:24", "08:30", "215A_LA PLATA", 113, "L6173"),
    @(21, "06:52:38", "08:34", "215A_LA PLATA", 102, "L6173"),
    @(22, "07:14:27", "08:35", "215A_LA PLATA", 81, "L6173"),
    @(23, "08:33:47", "08:37", "215A_LA PLATA", 4, "L6173"),
    @(24, "08:16:48", "09:08", "215D_LA PLATA", 52, "L6203"),
    @(25, "07:14:27", "09:09", "215D_LA PLATA", 115, "L6203"),
    @(26, "08:33:47", "09:10", "215D_LA PLATA", 37, "L6203"),
    @(27, "08:16:48", "10:02", "215B_LP-P MOR-40 Y 115", 106, "L6173")
)
foreach ($row in $rows) {
    $r = $row[0]
    if ($row[1] -ne $null) { $ws.Cells.Item($r, 1).Value = $row[1] }
    if ($row[2] -ne $null) { $ws.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 3).Value = $row[3] }
    if ($row[4] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[4] }
    if ($row[5] -ne $null) { $ws.Cells.Item($r, 5).Value = $row[5] }
}
